# Update the "Förändrad" (changed) date column C for rows 2-10
# from 2023-09-06 (serial 45175) to 2023-09-14 (serial 45183).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value = 45183
